$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: timestamp refined slightly (same displayed value, more precise fraction)
$ws.Range("A19").Value = 45874.75028924768

# New row 20 data
$ws.Range("A20").Value = 45874.79190705
$ws.Range("A20").NumberFormat = $ws.Range("A19").NumberFormat

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 16.39
$ws.Range("E20").Value = 83.69
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 2.87
$ws.Range("H20").Value = "E"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = "19:00:20"
